$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2-5 (practice rows): add carrier (column D) values
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Row 6-9: add pair_kind (column J) values
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: add kind (column C) and carrier (column D) values
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "look"
$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "look"
$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "where"
$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "where"
$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "can"
$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "can"
$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "do"
$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "do"
